$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 with new title
$ws.Range("C2").Value = "Права человека"

# Update H2 score
$ws.Range("H2").Value = 73

# Delete row 3 entirely (not admitted works removed)
$ws.Rows("3").Delete()
